$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.386.75"
$ws.Range("E2").Value = "  +3.67%  "

$ws.Range("D3").Value = "3.196.27"
$ws.Range("E3").Value = "  +5.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.87"
$ws.Range("E5").Value = "  +1.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "634.82"
$ws.Range("E6").Value = "  +0.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.235"
$ws.Range("E8").Value = "  +11.31%  "

$ws.Range("E9").Value = "  +5.64%  "

$ws.Range("D10").Value = "3.192.85"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.583"
$ws.Range("E11").Value = "  +33.31%  "

$ws.Range("E12").Value = "  +3.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.56"
$ws.Range("E13").Value = "  +8.98%  "

$ws.Range("D14").Value = "3.779.42"
$ws.Range("E14").Value = "  +4.96%  "

$ws.Range("E15").Value = "  +16.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.72"
$ws.Range("E16").Value = "  +7.38%  "

$ws.Range("D17").Value = "79.260.65"
$ws.Range("E17").Value = "  +3.46%  "

$ws.Range("D18").Value = "3.183.70"
$ws.Range("E18").Value = "  +4.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.58"
$ws.Range("E19").Value = "  +8.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.13"
$ws.Range("E20").Value = "  +35.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.14"
$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "429.05"
$ws.Range("E22").Value = "  +13.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.00"
$ws.Range("E23").Value = "  +14.33%  "

$ws.Range("B24").Value = "Aptos"
$ws.Range("C24").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.24"
$ws.Range("E24").Value = "  +12.65%  "

$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.358.84"
$ws.Range("E25").Value = "  +5.49%  "

$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.79"
$ws.Range("E26").Value = "  +8.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "76.73"
$ws.Range("E27").Value = "  +3.90%  "

$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000118"
$ws.Range("E29").Value = "  +5.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.99"
$ws.Range("E31").Value = "  +7.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.49"
$ws.Range("E32").Value = "  +4.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "527.59"
$ws.Range("E33").Value = "  +2.02%  "

$ws.Range("E34").Value = "  +1.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.143"
$ws.Range("E35").Value = "  +26.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.94"
$ws.Range("E36").Value = "  +9.73%  "

$ws.Range("E37").Value = "  +11.78%  "

$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.405"
$ws.Range("E39").Value = "  +5.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "164.61"
$ws.Range("E40").Value = "  +0.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.02"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "193.21"
$ws.Range("E42").Value = "  +2.62%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.47"
$ws.Range("E44").Value = "  +4.82%  "

$ws.Range("E45").Value = "  +10.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.81"
$ws.Range("E46").Value = "  +7.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.33"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.12"
$ws.Range("E48").Value = "  +2.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.91"
$ws.Range("E49").Value = "  +14.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.54"
$ws.Range("E50").Value = "  +3.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.633"
$ws.Range("E51").Value = "  +3.94%  "
